$wb = $excel.ActiveWorkbook

# --- "Light switch" sheet: tighten the force threshold ---------------------
$wsLightSwitch = $wb.Worksheets.Item("Light switch")
$wsLightSwitch.Range("C2").Value = ">15"

# --- "Typing" sheet: tighten the force threshold ----------------------------
$wsTyping = $wb.Worksheets.Item("Typing")
$wsTyping.Range("C2").Value = "<3"

# --- "Flight Mission Cycle" sheet: update the mission list ------------------
$wsFlightMission = $wb.Worksheets.Item("Flight Mission Cycle")
$wsFlightMission.Range("A2").Value = "Typing"
$wsFlightMission.Range("B2").Value = 1
$wsFlightMission.Range("A3").Value = "Piano"
$wsFlightMission.Range("B3").Value = 3
$wsFlightMission.Rows("4:4").Delete()

# --- "Piano" sheet: update the duration profile row -------------------------
$wsPiano = $wb.Worksheets.Item("Piano")
$wsPiano.Range("D3").Value = 10
$wsPiano.Range("E3").Value = 10
$wsPiano.Range("F3").Value = 5
$wsPiano.Range("G3").Value = 10
$wsPiano.Range("H3").Value = 5

# --- restore per-sheet selections to their new positions --------------------
$wsTyping.Activate()
$wsTyping.Range("H16").Select()

$wsLightSwitch.Activate()
$wsLightSwitch.Range("F9").Select()

$wsPiano.Activate()
$wsPiano.Range("G5").Select()

$wsFlightMission.Activate()
$wsFlightMission.Range("C10").Select()
